$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay stored as text (matches source formatting)
$textCells = @("D5","D6","D11","D15","D19","D20","D21","D24","D25","D30","D32","D33","D34","D35","D37","D38","D39","D41","D42","D43","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = "91.246.87"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.178.22"
$ws.Range("E3").Value = "  +5.14%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "216.45"
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("D6").Value = "628.50"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("E7").Value = "  +32.06%  "
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "3.173.67"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").Value = "0.771"
$ws.Range("E11").Value = "  +15.89%  "
$ws.Range("E12").Value = "  +8.18%  "
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("E14").Value = "  +7.05%  "
$ws.Range("D15").Value = "35.32"
$ws.Range("E15").Value = "  +10.35%  "
$ws.Range("D16").Value = "90.944.42"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "3.762.84"
$ws.Range("E17").Value = "  +4.92%  "
$ws.Range("D18").Value = "3.154.06"
$ws.Range("E18").Value = "  +3.73%  "
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  +13.15%  "
$ws.Range("D20").Value = "14.67"
$ws.Range("E20").Value = "  +9.97%  "
$ws.Range("D21").Value = "472.38"
$ws.Range("E21").Value = "  +11.90%  "
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("E23").Value = "  +11.94%  "
$ws.Range("D24").Value = "5.17"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "97.00"
$ws.Range("E25").Value = "  +18.01%  "
$ws.Range("E26").Value = "  +11.79%  "
$ws.Range("E27").Value = "  +8.58%  "
$ws.Range("D28").Value = "3.337.60"
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "9.39"
$ws.Range("E30").Value = "  +13.87%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "27.97"
$ws.Range("E33").Value = "  +23.60%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.193"
$ws.Range("E34").Value = "  +43.19%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "528.89"
$ws.Range("E35").Value = "  +6.10%  "
$ws.Range("E36").Value = "  +8.24%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.146"
$ws.Range("E37").Value = "  +10.63%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.65"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("D39").Value = "7.03"
$ws.Range("E39").Value = "  +6.38%  "
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("D41").Value = "0.0892"
$ws.Range("E41").Value = "  +29.09%  "
$ws.Range("D42").Value = "22.27"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "0.422"
$ws.Range("E43").Value = "  +18.10%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "2.00"
$ws.Range("E45").Value = "  +10.51%  "
$ws.Range("D47").Value = "0.714"
$ws.Range("E47").Value = "  +21.65%  "
$ws.Range("D48").Value = "152.26"
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("D49").Value = "4.68"
$ws.Range("E49").Value = "  +12.73%  "
$ws.Range("D50").Value = "1.37"
$ws.Range("E50").Value = "  +13.05%  "
$ws.Range("D51").Value = "45.35"
$ws.Range("E51").Value = "  +4.45%  "
